$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original values for rows 3-6 (these are the rows whose data gets
# rotated/reassigned per the source diff), for the columns that actually change.
$cols = @("A","B","D","E","F","G","H","Q","R")

$orig = @{}
foreach ($r in 3..6) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New row r gets the original values of row $srcMap[r]
$srcMap = @{
    3 = 4
    4 = 6
    5 = 3
    6 = 5
}

foreach ($r in 3..6) {
    $src = $srcMap[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $orig[$src][$c]
    }
}
